$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("week 13"): add homework hours and increase project hours
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 6.75

# Move the active selection to D14, matching the saved view state
$ws.Range("D14").Select()
